$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.760.91'
$ws.Range("E2").Value = '  -2.88%  '
$ws.Range("D3").Value = '3.329.90'
$ws.Range("E3").Value = '  -4.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '552.87'
$ws.Range("E5").Value = '  -2.98%  '
$ws.Range("D6").Value = '174.23'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("E7").Value = '  -2.95%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '3.323.37'
$ws.Range("E9").Value = '  -4.38%  '
$ws.Range("D10").Value = '0.619'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("E11").Value = '  +3.55%  '
$ws.Range("D12").Value = '53.85'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").Value = '8.99'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '3.865.79'
$ws.Range("E15").Value = '  -4.22%  '
$ws.Range("D16").Value = '18.14'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("E17").Value = '  -2.93%  '
$ws.Range("D18").Value = '3.331.90'
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("D19").Value = '11.79'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("D20").Value = '63.682.81'
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").Value = '0.969'
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").Value = '425.75'
$ws.Range("E22").Value = '  +3.62%  '
$ws.Range("D23").Value = '4.64'
$ws.Range("E23").Value = '  +10.43%  '
$ws.Range("D24").Value = '4.09'
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '83.82'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = '12.95'
$ws.Range("E26").Value = '  +2.87%  '
$ws.Range("D27").Value = '10.57'
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D28").Value = '2.80'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '8.66'
$ws.Range("E29").Value = '  -2.60%  '
$ws.Range("D30").Value = '29.49'
$ws.Range("E30").Value = '  -1.72%  '
$ws.Range("D31").Value = '6.55'
$ws.Range("E31").Value = '  +3.99%  '
$ws.Range("D32").Value = '587.68'
$ws.Range("E32").Value = '  -4.80%  '
$ws.Range("D33").Value = '11.38'
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("D35").Value = '58.46'
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  -5.11%  '
$ws.Range("D38").Value = '3.45'
$ws.Range("E38").Value = '  +3.65%  '
$ws.Range("D39").Value = '35.27'
$ws.Range("E39").Value = '  -4.07%  '
$ws.Range("D40").Value = '0.0₃0745'
$ws.Range("E40").Value = '  -5.28%  '
$ws.Range("E41").Value = '  -3.04%  '
$ws.Range("D42").Value = '3.096.62'
$ws.Range("E42").Value = '  -6.17%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").Value = '3.21'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").Value = '0.0404'
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("D47").Value = '2.42'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("E49").Value = '  -3.95%  '
$ws.Range("D50").Value = '133.03'
$ws.Range("E50").Value = '  -4.19%  '
$ws.Range("D51").Value = '8.15'
$ws.Range("E51").Value = '  -2.25%  '
